$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1: Summary
# ---------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B2").Value = 0.3327402135231317
$wsSummary.Range("C2").Value = 0.06516290726817042
$wsSummary.Range("D2").Value = 0.9285714285714286
$wsSummary.Range("E2").Value = 0.1217798594847775
$wsSummary.Range("F2").Value = 0.2544031311154599
$wsSummary.Range("G2").Value = 0.6151046405823476
$wsSummary.Range("H2").Value = 0.8103263777421081
$wsSummary.Range("I2").Value = 26
$wsSummary.Range("J2").Value = 373
$wsSummary.Range("K2").Value = 161
$wsSummary.Range("L2").Value = 2

# ---------------------------------------------------------------
# Sheet 2: Classification Report
# ---------------------------------------------------------------
$wsReport = $wb.Worksheets.Item("Classification Report")

# Row 2 ("0")
$wsReport.Range("B2").Value = 0.9877300613496932
$wsReport.Range("C2").Value = 0.301498127340824
$wsReport.Range("D2").Value = 0.4619799139167862

# Row 3 ("1")
$wsReport.Range("B3").Value = 0.06516290726817042
$wsReport.Range("C3").Value = 0.9285714285714286
$wsReport.Range("D3").Value = 0.1217798594847775

# Row 4 ("accuracy")
$wsReport.Range("B4").Value = 0.3327402135231317
$wsReport.Range("C4").Value = 0.3327402135231317
$wsReport.Range("D4").Value = 0.3327402135231317
$wsReport.Range("E4").Value = 0.3327402135231317

# Row 5 ("macro avg")
$wsReport.Range("B5").Value = 0.5264464843089318
$wsReport.Range("C5").Value = 0.6150347779561263
$wsReport.Range("D5").Value = 0.2918798867007819

# Row 6 ("weighted avg")
$wsReport.Range("B6").Value = 0.9417658615022153
$wsReport.Range("C6").Value = 0.3327402135231317
$wsReport.Range("D6").Value = 0.4450304450127003

# ---------------------------------------------------------------
# Sheet 3: Confusion Matrix
# ---------------------------------------------------------------
$wsMatrix = $wb.Worksheets.Item("Confusion Matrix")
$wsMatrix.Range("B2").Value = 161
$wsMatrix.Range("C2").Value = 373
$wsMatrix.Range("B3").Value = 2
$wsMatrix.Range("C3").Value = 26
